# manual dislocation uploading 2021/08/21 13:00
# Re-upload of the plan/fact dislocation dataset: the previously-entered
# August 2021 rows are superseded by a September 2021 extract pulled from
# the source system, and the cargo destination for the tail of the sheet
# changes from "Достык (эксп.)" to "Балхаш I".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: Row, ShippingDate (serial), CarAmount, FromStationName,
# ToStationName, CargoEtsngName
$rows = @(
    @(2, 44440, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(3, 44441, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(4, 44442, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(5, 44443, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(6, 44444, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(7, 44445, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(8, 44446, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(9, 44447, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(10, 44448, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(11, 44449, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(12, 44450, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(13, 44451, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(14, 44452, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(15, 44453, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(16, 44454, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(17, 44455, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(18, 44456, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(19, 44457, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(20, 44458, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(21, 44459, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(22, 44460, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(23, 44461, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(24, 44462, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(25, 44463, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(26, 44464, 0, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(27, 44465, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(28, 44466, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(29, 44467, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(30, 44468, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(31, 44469, 30, "Актогай", "Достык (эксп.)", "КОНЦЕНТР МЕД"),
    @(32, 44440, "s", "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(33, 44441, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(34, 44442, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(35, 44443, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(36, 44444, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(37, 44445, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(38, 44446, 20, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(39, 44447, 20, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(40, 44448, 15, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(41, 44449, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(42, 44450, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(43, 44451, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(44, 44452, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(45, 44453, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(46, 44454, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(47, 44455, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(48, 44456, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(49, 44457, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(50, 44458, 20, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(51, 44459, 20, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(52, 44460, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(53, 44461, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(54, 44462, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(55, 44463, 20, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(56, 44464, 15, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(57, 44465, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(58, 44466, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(59, 44467, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(60, 44468, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД"),
    @(61, 44469, 0, "Актогай", "Балхаш I", "КОНЦЕНТР МЕД")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Restore the view state: scrolled down to row 32, with G43 selected.
$ws.Range("G43").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 32
    $win.ScrollColumn = 1
} catch {
}
